$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.890.53"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.416.73"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.33"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.43"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E11").Value = "  -3.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.78"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "2.850.23"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "61.809.09"
$ws.Range("D17").Value = "2.414.23"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.28"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.43"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.79"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.47"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.73"
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "553.83"
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").Value = "2.536.91"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "0.0₃0930"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.15"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.146"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.80"
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.40"
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.51"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.62"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0525"
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.591"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.73"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0918"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  +0.71%  "
